# Insert a new weekly price record at row 303 in the "Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Tomate" sheet, pushing the existing rows
# 303:315 down to 304:316 (dimension grows from A1:R315 to A1:R316).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 303:315 down one row, leaving a blank (formatted) row 303.
$ws.Rows(303).Insert()

# Populate the newly inserted row 303 with the new record.
$ws.Range("A303").Value = 7
$ws.Range("B303").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C303").Value = "Ñuble"
$ws.Range("D303").Value = 44509
$ws.Range("E303").Value = 16
$ws.Range("F303").Value = 100112020
$ws.Range("G303").Value = "Tomate"
$ws.Range("H303").Value = "Larga vida"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 200
$ws.Range("K303").Value = 7500
$ws.Range("L303").Value = 8000
$ws.Range("M303").Value = 7750
$ws.Range("N303").Value = "$/caja 10 kilos"
$ws.Range("O303").Value = "Región de Arica y Parinacota"
$ws.Range("P303").Value = 775
$ws.Range("Q303").Value = 10
$ws.Range("R303").Value = "Hortaliza"
